$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate Dutch headers/descriptions into English
$ws.Range("A2").Value = "PP-number"
$ws.Range("B2").Value = "Group"
$ws.Range("C2").Value = "Ratio SDH I / II"

$ws.Range("A53").Value = "controls"
$ws.Range("A54").Value = "cyclists"

$ws.Range("A58").Value = "Groep: 1= cyclists"
$ws.Range("A59").Value = "Groep: 2=controls"

# Update the active selection to reflect where the author left off editing
$ws.Range("A60").Select()
